$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overview" (sheet1): add row 3 for the newly handed-back file
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3,1).Value = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsOverview.Cells.Item(3,2).Value = "'e2e\aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsOverview.Cells.Item(3,2).Style = "Hyperlink"
$wsOverview.Cells.Item(3,3).Value = ".md"
$wsOverview.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item(3,7).Value = "2016-10-18 12:11:29"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c10a184334c95831a813c5fd4eefd9020278aad/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md", "", "", "e2e\aa856427-b3d2-4f7e-a140-e7d1e57eef81.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ----------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): add row 3 for the newly handed-back file
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(3,1).Value = "'aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsZhCn.Cells.Item(3,1).Style = "Hyperlink"
$wsZhCn.Cells.Item(3,2).Value = ".md"
$wsZhCn.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$wsZhCn.Cells.Item(3,4).Value = "e2e"
$wsZhCn.Cells.Item(3,5).Value = "ht"
$wsZhCn.Cells.Item(3,6).Value = "'True"
$wsZhCn.Cells.Item(3,6).Style = "Normal"
$wsZhCn.Cells.Item(3,7).Value = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.zh-cn.xlf"
$wsZhCn.Cells.Item(3,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3,8).Value = "2016-10-18 12:11:16"
$wsZhCn.Cells.Item(3,9).Value = "'aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsZhCn.Cells.Item(3,9).Style = "Hyperlink"
$wsZhCn.Cells.Item(3,10).Value = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.zh-cn.xlf"
$wsZhCn.Cells.Item(3,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3,11).Value = "2016-10-18 12:11:58"
$wsZhCn.Cells.Item(3,12).Value = "'"
$wsZhCn.Cells.Item(3,12).Style = "Normal"
$wsZhCn.Cells.Item(3,13).Value = "'True"
$wsZhCn.Cells.Item(3,13).Style = "Normal"
$wsZhCn.Cells.Item(3,14).Value = "'"
$wsZhCn.Cells.Item(3,14).Style = "Normal"
$wsZhCn.Cells.Item(3,15).Value = "'False"
$wsZhCn.Cells.Item(3,15).Style = "Normal"
$wsZhCn.Cells.Item(3,16).Value = "'"
$wsZhCn.Cells.Item(3,16).Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/252c0608d6f1e5bf91db4b1bb704dd2022e2ba87/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md", "", "", "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/252c0608d6f1e5bf91db4b1bb704dd2022e2ba87/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md", "", "", "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ----------------------------------------------------------------------
# Sheet "de-de" (sheet3): add row 3 for the newly handed-back file
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(3,1).Value = "'aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsDeDe.Cells.Item(3,1).Style = "Hyperlink"
$wsDeDe.Cells.Item(3,2).Value = ".md"
$wsDeDe.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$wsDeDe.Cells.Item(3,4).Value = "e2e"
$wsDeDe.Cells.Item(3,5).Value = "ht"
$wsDeDe.Cells.Item(3,6).Value = "'True"
$wsDeDe.Cells.Item(3,6).Style = "Normal"
$wsDeDe.Cells.Item(3,7).Value = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.de-de.xlf"
$wsDeDe.Cells.Item(3,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3,8).Value = "2016-10-18 12:11:29"
$wsDeDe.Cells.Item(3,9).Value = "'aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$wsDeDe.Cells.Item(3,9).Style = "Hyperlink"
$wsDeDe.Cells.Item(3,10).Value = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.de-de.xlf"
$wsDeDe.Cells.Item(3,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3,11).Value = "2016-10-18 12:12:16"
$wsDeDe.Cells.Item(3,12).Value = "'"
$wsDeDe.Cells.Item(3,12).Style = "Normal"
$wsDeDe.Cells.Item(3,13).Value = "'True"
$wsDeDe.Cells.Item(3,13).Style = "Normal"
$wsDeDe.Cells.Item(3,14).Value = "'"
$wsDeDe.Cells.Item(3,14).Style = "Normal"
$wsDeDe.Cells.Item(3,15).Value = "'False"
$wsDeDe.Cells.Item(3,15).Style = "Normal"
$wsDeDe.Cells.Item(3,16).Value = "'"
$wsDeDe.Cells.Item(3,16).Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3425e9c2490bee39f574a571ff81a08c4f033aef/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md", "", "", "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3425e9c2490bee39f574a571ff81a08c4f033aef/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md", "", "", "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Host "Report generated for handback aa856427-b3d2-4f7e-a140-e7d1e57eef81"
